$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp shown in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 22:12"

# --- Reorder the block of country names in A169:A191 (rows for the ---
# --- countries that got shuffled in the source data refresh)       ---
$countries = @(
    "Guyana",
    "Angola",
    "Santa Sede",
    "Congo",
    "Suazilandia",
    "Guinea",
    "Mali",
    "Eritrea",
    "Cabo Verde",
    "Republica de Africa Central",
    "Republica del Chad",
    "Antigua y Barbuda",
    "San Martin (Parte Holandesa)",
    "Birmania",
    "Santa Lucia",
    "San Bartolome",
    "Mauritania",
    "Liberia",
    "Sudan",
    "Nepal",
    "Gambia",
    "Zimbabue",
    "Montserrat"
)

$startRow = 169
for ($i = 0; $i -lt $countries.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $countries[$i]
}

# --- Update the numeric statistics that changed for various rows ---
$ws.Cells.Item(4, 2).Value = 81996
$ws.Cells.Item(4, 3).Value = 13785
$ws.Cells.Item(4, 5).Value = 78955
$ws.Cells.Item(4, 7).Value = 150
$ws.Cells.Item(4, 8).Value = 1177

$ws.Cells.Item(7, 2).Value = 56347
$ws.Cells.Item(7, 3).Value = 6832
$ws.Cells.Item(7, 5).Value = 45178
$ws.Cells.Item(7, 7).Value = 507
$ws.Cells.Item(7, 8).Value = 4154

$ws.Cells.Item(17, 2).Value = 3910
$ws.Cells.Item(17, 3).Value = 501
$ws.Cells.Item(17, 5).Value = 3672
$ws.Cells.Item(17, 7).Value = 3
$ws.Cells.Item(17, 8).Value = 39

$ws.Cells.Item(28, 2).Value = 1819
$ws.Cells.Item(28, 3).Value = 255
$ws.Cells.Item(28, 5).Value = 1795
$ws.Cells.Item(28, 7).Value = 10
$ws.Cells.Item(28, 8).Value = 19

$ws.Cells.Item(54, 5).Value = 470
$ws.Cells.Item(54, 7).Value = 2
$ws.Cells.Item(54, 8).Value = 3

$ws.Cells.Item(103, 5).Value = 88
$ws.Cells.Item(103, 7).Value = 2
$ws.Cells.Item(103, 8).Value = 4

$ws.Cells.Item(170, 3).Value = 1

$ws.Cells.Item(173, 3).Value = 0

$ws.Cells.Item(175, 3).Value = 2

$ws.Cells.Item(176, 5).Value = 4
$ws.Cells.Item(176, 8).Value = 0

$ws.Cells.Item(177, 2).Value = 4
$ws.Cells.Item(177, 3).Value = 0
$ws.Cells.Item(177, 8).Value = 1

$ws.Cells.Item(185, 3).Value = 1
